# Sprint Backlog / Burndown workbook update
#
# The Week-1 column (D) previously had blank cells for the two newest
# backlog tasks ("Implement UI (fxml and codebehind)" in row 8 and
# "Prepare and make submission" in row 9). Fill those in with 0, matching
# the rest of the column, and leave the cursor on the next open row (D10)
# where the next task would be entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0

$ws.Range("D10").Select()
